# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
# to refreshed crypto data values, per commit "Updated symbol list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to stay text-typed (they were stored as literal
# text in the workbook, e.g. "315.24" / "3.10%", not numbers/percentages)
# while we overwrite their contents, then restore the default "Normal"
# style so no stray number-format/style is left attached to the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "315.24"
$ws.Cells.Item(2, 5).Value = "3.10%"
$ws.Cells.Item(3, 5).Value = "-0.03%"
$ws.Cells.Item(4, 4).Value = "5.088"
$ws.Cells.Item(4, 5).Value = "0.21%"
$ws.Cells.Item(5, 4).Value = "0.08158"
$ws.Cells.Item(5, 5).Value = "2.83%"
$ws.Cells.Item(6, 4).Value = "2.090"
$ws.Cells.Item(6, 5).Value = "-1.62%"
$ws.Cells.Item(7, 4).Value = "4.138"
$ws.Cells.Item(7, 5).Value = "0.41%"
$ws.Cells.Item(8, 4).Value = "7.939"
$ws.Cells.Item(8, 5).Value = "0.33%"
$ws.Cells.Item(9, 4).Value = "0.9345"
$ws.Cells.Item(9, 5).Value = "1.21%"
$ws.Cells.Item(10, 4).Value = "0.1039"
$ws.Cells.Item(10, 5).Value = "7.95%"
$ws.Cells.Item(11, 4).Value = "0.1919"
$ws.Cells.Item(11, 5).Value = "3.60%"
$ws.Cells.Item(12, 4).Value = "0.09068"
$ws.Cells.Item(12, 5).Value = "4.27%"
$ws.Cells.Item(13, 4).Value = "0.03600"
$ws.Cells.Item(13, 5).Value = "1.19%"
$ws.Cells.Item(14, 4).Value = "0.09875"
$ws.Cells.Item(14, 5).Value = "-0.41%"
$ws.Cells.Item(15, 4).Value = "0.001440"
$ws.Cells.Item(15, 5).Value = "-0.51%"
$ws.Cells.Item(16, 4).Value = "0.005842"
$ws.Cells.Item(17, 4).Value = "3.471"
$ws.Cells.Item(17, 5).Value = "-0.03%"
$ws.Cells.Item(18, 4).Value = "2.987"
$ws.Cells.Item(18, 5).Value = "8.51%"
$ws.Cells.Item(20, 4).Value = "0.1311"
$ws.Cells.Item(20, 5).Value = "-2.04%"
$ws.Cells.Item(21, 4).Value = "5.102"
$ws.Cells.Item(21, 5).Value = "-0.98%"
$ws.Cells.Item(22, 4).Value = "0.2215"
$ws.Cells.Item(22, 5).Value = "0.21%"
$ws.Cells.Item(23, 4).Value = "0.04555"
$ws.Cells.Item(23, 5).Value = "0.83%"
$ws.Cells.Item(24, 4).Value = "0.001242"
$ws.Cells.Item(24, 5).Value = "0.70%"
$ws.Cells.Item(25, 4).Value = "0.004801"
$ws.Cells.Item(25, 5).Value = "-1.14%"
$ws.Cells.Item(26, 5).Value = "-3.50%"
$ws.Cells.Item(27, 4).Value = "0.0004506"
$ws.Cells.Item(27, 5).Value = "-5.24%"
$ws.Cells.Item(39, 4).Value = "0.01967"
$ws.Cells.Item(39, 5).Value = "6.46%"
$ws.Cells.Item(40, 4).Value = "0.04920"
$ws.Cells.Item(40, 5).Value = "3.53%"
$ws.Cells.Item(41, 4).Value = "0.007599"
$ws.Cells.Item(41, 5).Value = "-2.61%"
$ws.Cells.Item(42, 4).Value = "0.1383"
$ws.Cells.Item(42, 5).Value = "-1.37%"
$ws.Cells.Item(43, 4).Value = "0.007879"
$ws.Cells.Item(43, 5).Value = "1.56%"
$ws.Cells.Item(44, 5).Value = "-3.97%"
$ws.Cells.Item(45, 4).Value = "0.01173"
$ws.Cells.Item(45, 5).Value = "5.03%"
$ws.Cells.Item(46, 4).Value = "0.00006718"
$ws.Cells.Item(46, 5).Value = "6.82%"
$ws.Cells.Item(47, 4).Value = "0.00000000751"
$ws.Cells.Item(47, 5).Value = "0.40%"
$ws.Cells.Item(48, 4).Value = "156.88"
$ws.Cells.Item(48, 5).Value = "209.72%"
$ws.Cells.Item(49, 4).Value = "0.001702"
$ws.Cells.Item(49, 5).Value = "-10.52%"
$ws.Cells.Item(50, 4).Value = "0.00002103"
$ws.Cells.Item(50, 5).Value = "0.40%"
$ws.Cells.Item(51, 4).Value = "0.0002003"
$ws.Cells.Item(51, 5).Value = "0.40%"

# Restore default styling on the touched range (remove the temporary text
# number format so cells look exactly as they did before, just with new text).
$dataRange.Style = "Normal"

